{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the target paragraphs by their (unique) text content rather than\n// relying on fixed indices, so the script is resilient to minor shifts.\nlet pIngresos = null;       // \"Opci\u00f3n de Ingresos y egresos varios\"\nlet pCamposSig = null;      // \"\u2026campos siguientes\"\nlet pServicioDom = null;    // \"En Servicio a Domicilio...\"\nlet pDescFactura = null;    // \"Opci\u00f3n de ... Descripci\u00f3n de Factura ... personalizado (Ej. Por consumo)\"\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"Opci\u00f3n de Ingresos y egresos varios\") !== -1) {\n    pIngresos = paras.items[i];\n  } else if (t.indexOf(\"\u2026campos siguientes\") !== -1) {\n    pCamposSig = paras.items[i];\n  } else if (t.indexOf(\"En Servicio a Domicilio\") !== -1) {\n    pServicioDom = paras.items[i];\n  } else if (t.indexOf(\"Descripci\u00f3n de Factura\") !== -1) {\n    pDescFactura = paras.items[i];\n  }\n}\n\n// 1) \"Opci\u00f3n de Ingresos y egresos varios\" -> make it green (00B050),\n//    applied both to the paragraph mark and the run.\nif (pIngresos) {\n  pIngresos.font.color = \"#00B050\";\n}\n\n// 3) \"En Servicio a Domicilio...\" -> change color from dark blue (002060)\n//    to green (00B050).\nif (pServicioDom) {\n  pServicioDom.font.color = \"#00B050\";\n}\n\nawait context.sync();\n\n// 4) Merge the \" persona\" + \"lizado (Ej. Por consumo)\" runs into a single\n//    run reading \" personalizado (Ej. Por consumo)\", and drop the\n//    \"_GoBack\" bookmark that previously sat between them.\nif (pDescFactura) {\n  const personaResults = pDescFactura.search(\" persona\", { matchCase: true });\n  personaResults.load(\"items\");\n  await context.sync();\n\n  if (personaResults.items.length > 0) {\n    const personaRange = personaResults.items[0];\n    const boundaryPoint = personaRange.getRange(\"End\");\n    const paraEnd = pDescFactura.getRange(\"End\");\n    const tailRange = boundaryPoint.expandTo(paraEnd);\n\n    // Remove the trailing \"lizado (Ej. Por consumo)\" run (and the bookmark\n    // that sits right before it).\n    tailRange.delete();\n    await context.sync();\n\n    // Replace \" persona\" with the merged text \" personalizado (Ej. Por consumo)\".\n    personaRange.insertText(\" personalizado (Ej. Por consumo)\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2) \"\u2026campos siguientes\" -> insert the \"_GoBack\" bookmark right before it\n//    (moved from its old location inside the invoice-description paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nif (pCamposSig) {\n  const startRange = pCamposSig.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Green color used throughout this document (hex 00B050 -> decimal BGR-ish\n# long value used by Word's Font.Color: R + G*256 + B*65536).\n$green = 5287936\n\n# ---------------------------------------------------------------------\n# 1) \"Opcion de Ingresos y egresos varios\" -> apply green font color.\n# ---------------------------------------------------------------------\n$pIngresos = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*Ingresos y egresos varios*\") {\n    $pIngresos = $p\n    break\n  }\n}\nif ($pIngresos -ne $null) {\n  $pIngresos.Range.Font.Color = $green\n}\n\n# ---------------------------------------------------------------------\n# 3) \"En Servicio a Domicilio...\" -> change color from dark blue (002060)\n#    to green (00B050).\n# ---------------------------------------------------------------------\n$pServicio = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*En Servicio a Domicilio*\") {\n    $pServicio = $p\n    break\n  }\n}\nif ($pServicio -ne $null) {\n  $pServicio.Range.Font.Color = $green\n}\n\n# ---------------------------------------------------------------------\n# 4) Merge the \" persona\" + \"lizado (Ej. Por consumo)\" runs into a single\n#    run reading \" personalizado (Ej. Por consumo)\". This also removes the\n#    \"_GoBack\" bookmark that used to sit between them (it gets re-added in\n#    step 2 below, right before \"...campos siguientes\").\n# ---------------------------------------------------------------------\n$pDescFactura = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if (($p.Range.Text -like \"*Factura*\") -and ($p.Range.Text -like \"*consumo*\")) {\n    $pDescFactura = $p\n    break\n  }\n}\nif ($pDescFactura -ne $null) {\n  $searchRange = $pDescFactura.Range.Duplicate\n  $found = $searchRange.Find.Execute(\" persona\")\n  if ($found) {\n    $startPos = $searchRange.Start\n    $endPos = $searchRange.End\n\n    # Remove the trailing \"lizado (Ej. Por consumo)\" text (everything after\n    # \" persona\" up to, but excluding, the paragraph mark).\n    $paraEnd = $pDescFactura.Range.End\n    $tailRange = $d.Range($endPos, $paraEnd - 1)\n    $tailRange.Delete()\n\n    # Replace \" persona\" with \" personalizado (Ej. Por consumo)\". Toggling\n    # Bold on/off around the text assignment forces Word to keep this as\n    # its own run instead of silently re-merging it into the preceding\n    # \"Descripcion de Factura\" run (which has identical formatting).\n    $freshRange = $d.Range($startPos, $startPos + (\" persona\").Length)\n    $freshRange.Bold = 1\n    $freshRange.Text = \" personalizado (Ej. Por consumo)\"\n    $freshRange2 = $d.Range($startPos, $startPos + (\" personalizado (Ej. Por consumo)\").Length)\n    $freshRange2.Bold = 0\n  }\n}\n\n# ---------------------------------------------------------------------\n# 2) Move the \"_GoBack\" bookmark so it sits right before \"...campos\n#    siguientes\" instead of inside the invoice-description paragraph.\n# ---------------------------------------------------------------------\ntry {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n} catch {\n}\n\n$pCamposSig = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*campos siguientes*\") {\n    $pCamposSig = $p\n    break\n  }\n}\nif ($pCamposSig -ne $null) {\n  $bmRange = $pCamposSig.Range.Duplicate\n  $bmRange.Collapse(1)\n  $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
